# Apply scheduled-runner data updates to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2177.16
$ws.Range("I137").Value = 2779.3076
$ws.Range("J137").Value = 1524.8334
$ws.Range("K137").Value = 8337.9228
$ws.Range("L137").Value = 4574.5002
$ws.Range("M137").Value = -5787.9228
$ws.Range("N137").Value = -9674.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 347282.53
$ws.Range("I61").Value = 429563.72
$ws.Range("J61").Value = 1701.5
$ws.Range("K61").Value = 429563.72
$ws.Range("L61").Value = 1701.5
$ws.Range("M61").Value = -429351.72
$ws.Range("N61").Value = -2125.5
$ws.Range("H74").Value = 20001576
$ws.Range("I74").Value = 21740636
$ws.Range("J74").Value = 2378
$ws.Range("K74").Value = 21740636
$ws.Range("L74").Value = 2378
$ws.Range("M74").Value = -21739762
$ws.Range("N74").Value = -4126
$ws.Range("H77").Value = 20001576
$ws.Range("I77").Value = 21740636
$ws.Range("J77").Value = 2378
$ws.Range("K77").Value = 108703180
$ws.Range("L77").Value = 11890
$ws.Range("M77").Value = -108698812
$ws.Range("N77").Value = -20626
$ws.Range("H108").Value = 42842.5
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 42842.5
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 42842.5
$ws.Range("N108").Value = -50522.5
$ws.Range("H132").Value = 8187.3716
$ws.Range("I132").Value = 1019.79663
$ws.Range("J132").Value = 46631.637
$ws.Range("K132").Value = 3059.38989
$ws.Range("L132").Value = 139894.911
$ws.Range("M132").Value = -529.3898900000004
$ws.Range("N132").Value = -144954.911
$ws.Range("H136").Value = 347282.53
$ws.Range("I136").Value = 429563.72
$ws.Range("J136").Value = 1701.5
$ws.Range("K136").Value = 1288691.16
$ws.Range("L136").Value = 5104.5
$ws.Range("M136").Value = -1286141.16
$ws.Range("N136").Value = -10204.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2859
$ws.Range("I134").Value = 3024.3242
$ws.Range("J134").Value = 2179.3333
$ws.Range("K134").Value = 9072.972600000001
$ws.Range("L134").Value = 6537.999899999999
$ws.Range("M134").Value = -6537.972600000001
$ws.Range("N134").Value = -11607.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4800.2085
$ws.Range("I31").Value = 3323.6365
$ws.Range("J31").Value = 6049.615
$ws.Range("K31").Value = 3323.6365
$ws.Range("L31").Value = 6049.615
$ws.Range("M31").Value = -3028.6365
$ws.Range("N31").Value = -6639.615
$ws.Range("H34").Value = 4800.2085
$ws.Range("I34").Value = 3323.6365
$ws.Range("J34").Value = 6049.615
$ws.Range("K34").Value = 3323.6365
$ws.Range("L34").Value = 6049.615
$ws.Range("M34").Value = -3121.6365
$ws.Range("N34").Value = -6453.615
$ws.Range("H132").Value = 1599.1296
$ws.Range("I132").Value = 1292.3
$ws.Range("J132").Value = 5434.5
$ws.Range("K132").Value = 3876.9
$ws.Range("L132").Value = 16303.5
$ws.Range("M132").Value = -1346.9
$ws.Range("N132").Value = -21363.5
$ws.Range("H134").Value = 648.26086
$ws.Range("I134").Value = 692.56665
$ws.Range("J134").Value = 565.1875
$ws.Range("K134").Value = 2077.69995
$ws.Range("L134").Value = 1695.5625
$ws.Range("M134").Value = 457.3000499999998
$ws.Range("N134").Value = -6765.5625
$ws.Range("H137").Value = 25000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 25000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 25000
$ws.Range("N137").Value = -35200
$ws.Range("H140").Value = 50780
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 50780
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 50780
$ws.Range("N140").Value = -61140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 1994
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 1994
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 5982
$ws.Range("N48").Value = -6482
$ws.Range("H68").Value = 933
$ws.Range("I68").Value = 900
$ws.Range("J68").Value = 999
$ws.Range("K68").Value = 2700
$ws.Range("L68").Value = 2997
$ws.Range("M68").Value = -1889
$ws.Range("N68").Value = -4619
$ws.Range("H71").Value = 933
$ws.Range("I71").Value = 900
$ws.Range("J71").Value = 999
$ws.Range("K71").Value = 8100
$ws.Range("L71").Value = 8991
$ws.Range("M71").Value = -4044
$ws.Range("N71").Value = -17103
$ws.Range("H92").Value = 676.6667
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 765
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 2295
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -4791
$ws.Range("H131").Value = 850.15
$ws.Range("I131").Value = 550
$ws.Range("J131").Value = 862.65625
$ws.Range("K131").Value = 1650
$ws.Range("L131").Value = 2587.96875
$ws.Range("M131").Value = 3390
$ws.Range("N131").Value = -12667.96875
$ws.Range("H132").Value = 679.8
$ws.Range("H137").Value = 18521526
$ws.Range("I137").Value = 1275.8
$ws.Range("J137").Value = 25644700
$ws.Range("K137").Value = 3827.4
$ws.Range("L137").Value = 76934100
$ws.Range("M137").Value = 1272.6
$ws.Range("N137").Value = -76944300

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 27400
$ws.Range("I57").Value = 27400
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 27400
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -26580
$ws.Range("N57").ClearContents()
$ws.Range("H80").Value = 3586.3635
$ws.Range("I80").Value = 2887.5
$ws.Range("J80").Value = 3985.7144
$ws.Range("K80").Value = 2887.5
$ws.Range("L80").Value = 3985.7144
$ws.Range("M80").Value = -1889.5
$ws.Range("N80").Value = -5981.7144
$ws.Range("H83").Value = 3586.3635
$ws.Range("I83").Value = 2887.5
$ws.Range("J83").Value = 3985.7144
$ws.Range("K83").Value = 14437.5
$ws.Range("L83").Value = 19928.572
$ws.Range("M83").Value = -9445.5
$ws.Range("N83").Value = -29912.572
$ws.Range("H132").Value = 18571.033
$ws.Range("I132").Value = 2541.4482
$ws.Range("J132").Value = 251000
$ws.Range("K132").Value = 7624.344599999999
$ws.Range("L132").Value = 753000
$ws.Range("M132").Value = -5094.344599999999
$ws.Range("N132").Value = -758060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3518.0588
$ws.Range("I40").Value = 3154
$ws.Range("J40").Value = 4701.25
$ws.Range("K40").Value = 3154
$ws.Range("L40").Value = 4701.25
$ws.Range("M40").Value = -3018
$ws.Range("N40").Value = -4973.25
$ws.Range("H132").Value = 1210.0731
$ws.Range("I132").Value = 1210.0731
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3630.2193
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1100.2193
$ws.Range("H136").Value = 972.0789
$ws.Range("I136").Value = 910
$ws.Range("J136").Value = 1381.8
$ws.Range("K136").Value = 2730
$ws.Range("L136").Value = 4145.4
$ws.Range("M136").Value = -180
$ws.Range("N136").Value = -9245.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H136").Value = 19233074
$ws.Range("I136").Value = 27028078
$ws.Range("J136").Value = 5400.933
$ws.Range("K136").Value = 81084234
$ws.Range("L136").Value = 16202.799
$ws.Range("M136").Value = -81081684
$ws.Range("N136").Value = -21302.799
